$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-12 (word, total frequency, websites)
$data = @(
    @("Cookies",   83,  "https://www.100-dakar.com (14), https://www.benjaminwahl.at (18), https://www.drehorgelkabarett.at (14), https://www.ottosaxinger.at (3), https://www.peligro.at (14), https://www.schuledesungehorsams.at (2), https://www.skodone.at (18)"),
    @("Page",      35,  "https://www.hungaromedia.at (8), https://www.kuenstlerinnen.at (8), https://www.luckeneder-art.at (8), https://www.platform-socialism.org (3), https://www.regional-express.org (8)"),
    @("Linz",      52,  "https://www.freie-medien.at (10), https://www.freizeitundkommunikation.at (3), https://www.linzfmr.at (18), https://www.steingeschichten.at (21)"),
    @("März",      49,  "https://www.das-kollektiv.at (8), https://www.feminismus-krawall.at (16), https://www.fiftitu.at (19), https://www.unkraut-comics.at (6)"),
    @("Art",       29,  "https://www.eipcp.net (13), https://www.kairus.org (6), https://www.negentropy-sport.net (2), https://www.radical-openness.org (8)"),
    @("Schule",    294, "https://www.alteschule-gutau.at (4), https://www.derschueler.at (5), https://www.die-schule.at (285)"),
    @("Uhr",       35,  "https://www.fro.at (23), https://www.rudolfhabringer.at (12)"),
    @("Kultur",    31,  "https://www.frauenkultur.at (17), https://www.igkultur.at (14)"),
    @("Film",      15,  "https://www.corpushomini.info (3), https://www.doublehappiness.at (12)"),
    @("Magdalena", 12,  "https://www.magdalenareiter.at (2), https://www.themagdalenaproject.org (10)"),
    @("Andreas",   8,   "https://www.andreaskurz.net (2), https://www.andreaszingerle.com (6)")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Remove the leftover rows (previously 13-16) that no longer exist in the new data
$ws.Range("A13:C16").Clear()
